$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 86.31999969482422
$ws.Cells.Item(2, 3).Value = 81.62000274658203
$ws.Cells.Item(2, 4).Value = -1.450000047683716
$ws.Cells.Item(2, 7).Value = 5.639999866485596
$ws.Cells.Item(2, 8).Value = 4
$ws.Cells.Item(2, 9).Value = 3.849999904632568
$ws.Cells.Item(2, 10).Value = 3.34
$ws.Cells.Item(2, 13).Value = 26.36000061035156
$ws.Cells.Item(2, 14).Value = 10.67
$ws.Cells.Item(2, 15).Value = 2.5
$ws.Cells.Item(2, 16).Value = 2.42
$ws.Cells.Item(2, 17).Value = -322.25
$ws.Cells.Item(2, 18).Value = 72.54000000000001
$ws.Cells.Item(3, 2).Value = 22.80999946594238
$ws.Cells.Item(3, 3).Value = 24.72999954223633
$ws.Cells.Item(3, 4).Value = -1.639999985694885
$ws.Cells.Item(3, 5).Value = 53.11999893188477
$ws.Cells.Item(3, 6).Value = 7.73
$ws.Cells.Item(3, 7).Value = 4.110000133514404
$ws.Cells.Item(3, 8).Value = 3.8
$ws.Cells.Item(3, 9).Value = 2.210000038146973
$ws.Cells.Item(3, 10).Value = 2.39
$ws.Cells.Item(3, 11).Value = 6.230000019073486
$ws.Cells.Item(3, 12).Value = 3.12
$ws.Cells.Item(3, 13).Value = 32.04999923706055
$ws.Cells.Item(3, 14).Value = 5.81
$ws.Cells.Item(3, 15).Value = 1.25
$ws.Cells.Item(3, 16).Value = 1.81
$ws.Cells.Item(3, 17).Value = -338.3099975585938
$ws.Cells.Item(3, 18).Value = 74.48999999999999
$ws.Cells.Item(4, 2).Value = 97.87999725341797
$ws.Cells.Item(4, 3).Value = 73.70999908447266
$ws.Cells.Item(4, 4).Value = -1.450000047683716
$ws.Cells.Item(4, 5).Value = 51.54999923706055
$ws.Cells.Item(4, 6).Value = 8.44
$ws.Cells.Item(4, 7).Value = 5.650000095367432
$ws.Cells.Item(4, 8).Value = 4
$ws.Cells.Item(4, 9).Value = 3.849999904632568
$ws.Cells.Item(4, 10).Value = 3.34
$ws.Cells.Item(4, 11).Value = 7.880000114440918
$ws.Cells.Item(4, 12).Value = 4.25
$ws.Cells.Item(4, 13).Value = 26.27000045776367
$ws.Cells.Item(4, 14).Value = 10.67
$ws.Cells.Item(4, 15).Value = 2.5
$ws.Cells.Item(4, 16).Value = 2.42
$ws.Cells.Item(4, 17).Value = -322.1799926757812
$ws.Cells.Item(4, 18).Value = 72.48999999999999
$ws.Cells.Item(5, 2).Value = 30.48999977111816
$ws.Cells.Item(5, 3).Value = 24.94000053405762
$ws.Cells.Item(5, 4).Value = -1.549999952316284
$ws.Cells.Item(5, 5).Value = 48.91999816894531
$ws.Cells.Item(5, 6).Value = 7.99
$ws.Cells.Item(5, 7).Value = 3.609999895095825
$ws.Cells.Item(5, 8).Value = 3.9
$ws.Cells.Item(5, 9).Value = 0.5400000214576721
$ws.Cells.Item(5, 10).Value = 2.12
$ws.Cells.Item(5, 11).Value = 3.390000104904175
$ws.Cells.Item(5, 12).Value = 3.05
$ws.Cells.Item(5, 13).Value = 49.27999877929688
$ws.Cells.Item(5, 14).Value = 6.09
$ws.Cells.Item(5, 15).Value = 0.6399999856948853
$ws.Cells.Item(5, 16).Value = 0.83
$ws.Cells.Item(5, 17).Value = -356.4200134277344
$ws.Cells.Item(5, 18).Value = 73.59
$ws.Cells.Item(6, 2).Value = 80.81999969482422
$ws.Cells.Item(6, 3).Value = 14.15999984741211
$ws.Cells.Item(6, 4).Value = -1.580000042915344
$ws.Cells.Item(6, 5).Value = 51.61999893188477
$ws.Cells.Item(6, 6).Value = 8.44
$ws.Cells.Item(6, 7).Value = 5.269999980926514
$ws.Cells.Item(6, 8).Value = 3.98
$ws.Cells.Item(6, 9).Value = 3.849999904632568
$ws.Cells.Item(6, 10).Value = 3.34
$ws.Cells.Item(6, 11).Value = 7.869999885559082
$ws.Cells.Item(6, 12).Value = 4.25
$ws.Cells.Item(6, 13).Value = 27.03000068664551
$ws.Cells.Item(6, 14).Value = 10.62
$ws.Cells.Item(6, 15).Value = 2.5
$ws.Cells.Item(6, 16).Value = 2.42
$ws.Cells.Item(6, 17).Value = -323.5
$ws.Cells.Item(6, 18).Value = 73.25
